# Update the raw experimental data in "Planilha1" (B17:D21).
# Columns F,G,I,J,L,M on these rows (and the dependent rows 56-67) are
# formulas that recompute automatically, as do the three embedded charts
# that plot this range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# B = substrate concentration [S], C = reaction velocity v, D = second
# experiment's velocity (now reset to 0 for all rows).
$ws.Range("B17").Value = 20
$ws.Range("C17").Value = 0.32
$ws.Range("D17").Value = 0

$ws.Range("B18").Value = 50
$ws.Range("C18").Value = 0.521
$ws.Range("D18").Value = 0

$ws.Range("B19").Value = 100
$ws.Range("C19").Value = 0.699
$ws.Range("D19").Value = 0

$ws.Range("B20").Value = 300
$ws.Range("C20").Value = 1.11
$ws.Range("D20").Value = 0

$ws.Range("B21").Value = 700
$ws.Range("C21").Value = 1.33
$ws.Range("D21").Value = 0

$excel.Calculate()

# Reflect the author's updated view state (zoom + selection) as closely as
# the host lets us.
$excel.ActiveWindow.Zoom = 107
$ws.Range("K29").Select() | Out-Null
